$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.383.70'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.420.29'
$ws.Range('E3').Value = '  -0.91%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '584.50'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.21'
$ws.Range('E6').Value = '  +1.52%  '
$ws.Range('E7').Value = '  +3.70%  '
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.421.14'
$ws.Range('E9').Value = '  -0.76%  '
$ws.Range('E10').Value = '  +0.25%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.95'
$ws.Range('E11').Value = '  +1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.415'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.010.39'
$ws.Range('E13').Value = '  -1.51%  '
$ws.Range('E14').Value = '  +0.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.50'
$ws.Range('E15').Value = '  -2.60%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.412.13'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.415.76'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.92'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.80'
$ws.Range('E20').Value = '  -0.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '368.99'
$ws.Range('E21').Value = '  -2.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.61'
$ws.Range('E22').Value = '  -2.58%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.84'
$ws.Range('E23').Value = '  +0.68%  '
$ws.Range('E24').Value = '  +5.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.996'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('E26').Value = '  -0.50%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.85'
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.82'
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '23.49'
$ws.Range('E32').Value = '  -3.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.05'
$ws.Range('E34').Value = '  -1.77%  '
$ws.Range('E35').Value = '  -3.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.55'
$ws.Range('E36').Value = '  -1.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.95'
$ws.Range('E37').Value = '  +1.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.877'
$ws.Range('E38').Value = '  -1.66%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.71'
$ws.Range('E39').Value = '  -5.94%  '
$ws.Range('E40').Value = '  +0.66%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.61'
$ws.Range('E41').Value = '  -1.60%  '
$ws.Range('E42').Value = '  -1.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.713.77'
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.31'
$ws.Range('E44').Value = '  -3.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0691'
$ws.Range('E45').Value = '  -0.59%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.35'
$ws.Range('E46').Value = '  +2.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.11'
$ws.Range('E47').Value = '  -1.37%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '334.19'
$ws.Range('E48').Value = '  +7.90%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0286'
$ws.Range('E49').Value = '  -2.80%  '
$ws.Range('E50').Value = '  +2.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.90'
$ws.Range('E51').Value = '  +3.88%  '
